# Adding modules for real-time LPF control, fixed Butterworth filter
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Activate()

# G5: 16 -> 32
$ws.Range("G5").Value = 32

# I5: "200 kHz" -> "250 kHz"
$ws.Range("I5").Value = "250 kHz"

# Update the view/scroll position (best effort) and selection to match
# the saved workbook state (topLeftCell J2, active cell I5).
$excel.ActiveWindow.ScrollColumn = 10
$excel.ActiveWindow.ScrollRow = 2
$ws.Range("I5").Select()
